# Atualização de bases das ligas, do dia: 17-05-2024 às 13:59
# Swap the data (all columns except A = row index) between row pairs
# 36/37 and 122/123 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $row1, $row2) {
    # Columns B..AB (2..28) hold the match data; column A is the row index
    # and must stay untouched.
    $firstCol = 2
    $lastCol = 28

    $range1 = $ws.Range($ws.Cells.Item($row1, $firstCol), $ws.Cells.Item($row1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($row2, $firstCol), $ws.Cells.Item($row2, $lastCol))

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

Swap-Rows $ws 36 37
Swap-Rows $ws 122 123
